# Fix bug in VBA code; see #805
# Mirror the computed "Sample calcs" results (column C) into a new plain-value
# column D, add the missing "K" unit label next to the dewpoint result, and
# leave the sheet selected/sized the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample calcs")

# Mirror the cached formula results from column C into column D as plain values
$ws.Range("D6").Value2  = $ws.Range("C6").Value2
$ws.Range("D7").Value2  = $ws.Range("C7").Value2
$ws.Range("D8").Value2  = $ws.Range("C8").Value2
$ws.Range("D9").Value2  = $ws.Range("C9").Value2
$ws.Range("D10").Value2 = $ws.Range("C10").Value2
$ws.Range("D12").Value2 = $ws.Range("C12").Value2

# Dewpoint-of-dry-air row: add the missing "K" unit label and mirror the result
$ws.Range("B20").Value2 = "K"
$ws.Range("D20").Value2 = $ws.Range("C20").Value2

# Narrow column C back down now that column D carries the duplicate values
$ws.Columns.Item(3).ColumnWidth = 14.7

# Leave the selection where the author finished editing
$ws.Range("C20").Select() | Out-Null
